$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shared-string (rich text) edits
#    A8  = "Volume 31   Number  31"  -> "Volume 31   Number  32"
#    C9  = "Report Covering the Week  7/29/2024  Through  8/4/2024"
#          -> "Report Covering the Week  8/5/2024  Through  8/11/2024"
# ---------------------------------------------------------------------------

# A8: second "31" (after "Number") -> "32"
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "32"

# C9: replace the second date first so the first date's character offsets
# are not disturbed by a length change.
$c9 = $ws.Range("C9")
$c9.Characters(47, 8).Text = "8/11/2024"
$c9.Characters(27, 9).Text = "8/5/2024"

# ---------------------------------------------------------------------------
# 2. Cells that change data type (number <-> text) in rows 23, 27, 28.
#    These need both the value AND the original cell style restored, since
#    assigning a text value that looks numeric nudges Excel onto a
#    "quote prefix" style. We fix that by pasting the number *format* back
#    in from an untouched donor cell with the desired style.
# ---------------------------------------------------------------------------

# C23: number 2 -> text "0" (style 14, same as D14/C14/F14/G14)
$ws.Range("C23").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null

# D23: text "0" -> number 1 (style 15, same as I14/J14)
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Value = 1

# E23: text "***.*" -> number -100 (style 16, same as K14/L14/M14)
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = -100

# D27: number 1 -> text "0" (style 14)
$ws.Range("D27").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null

# E27: number 100 -> text "***.*" (style 14, same as E14)
$ws.Range("E27").Value = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null

# C28: text "0" -> number 2 (style 15)
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value = 2

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Plain numeric value updates (rows 14-30), style/type unchanged.
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -66.666666666666
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 4
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = -38.461538461538
$ws.Range("M15").Value = 33.333333333333
$ws.Range("N15").Value = -27.272727272727
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 96
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -17.241379310344
$ws.Range("M16").Value = -1.030927835051
$ws.Range("N16").Value = -73.913043478260
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -13.636363636363
$ws.Range("I17").Value = 149
$ws.Range("J17").Value = 142
$ws.Range("K17").Value = 4.929577464788
$ws.Range("L17").Value = 7.971014492753
$ws.Range("M17").Value = 81.707317073170
$ws.Range("N17").Value = -11.309523809523
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 228.571428571429
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 97
$ws.Range("K18").Value = -3.092783505154
$ws.Range("L18").Value = 13.253012048192
$ws.Range("M18").Value = 10.588235294117
$ws.Range("N18").Value = -86.416184971098
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 88.888888888888
$ws.Range("F19").Value = 86
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = 91.111111111111
$ws.Range("I19").Value = 496
$ws.Range("J19").Value = 406
$ws.Range("K19").Value = 22.167487684729
$ws.Range("L19").Value = 19.806763285024
$ws.Range("M19").Value = 130.697674418605
$ws.Range("N19").Value = 90.769230769230
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -13.793103448275
$ws.Range("I20").Value = 204
$ws.Range("J20").Value = 208
$ws.Range("K20").Value = -1.923076923076
$ws.Range("L20").Value = 23.636363636363
$ws.Range("M20").Value = 164.935064935065
$ws.Range("N20").Value = -82.073813708260
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 34.782608695652
$ws.Range("F21").Value = 172
$ws.Range("G21").Value = 115
$ws.Range("H21").Value = 49.565217391304
$ws.Range("I21").Value = 1048
$ws.Range("J21").Value = 959
$ws.Range("K21").Value = 9.280500521376
$ws.Range("L21").Value = 12.567132116004
$ws.Range("M21").Value = 85.159010600706
$ws.Range("N21").Value = -60.303030303030
$ws.Range("M22").Value = -78.571428571428
$ws.Range("J23").Value = 54
$ws.Range("K23").Value = -20.370370370370
$ws.Range("M23").Value = 53.571428571428
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -68.965517241379
$ws.Range("F24").Value = 75
$ws.Range("H24").Value = -21.875
$ws.Range("I24").Value = 678
$ws.Range("J24").Value = 668
$ws.Range("K24").Value = 1.497005988023
$ws.Range("L24").Value = -7.503410641200
$ws.Range("M24").Value = 39.506172839506
$ws.Range("C25").Value = 9
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = 4
$ws.Range("I25").Value = 272
$ws.Range("J25").Value = 267
$ws.Range("K25").Value = 1.872659176029
$ws.Range("L25").Value = -32.673267326732
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -77.777777777777
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 205
$ws.Range("J26").Value = 194
$ws.Range("K26").Value = 5.670103092783
$ws.Range("L26").Value = -3.301886792452
$ws.Range("M26").Value = -14.937759336099
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 16
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -11.111111111111
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("J28").Value = 19
$ws.Range("K28").Value = 89.473684210526
$ws.Range("N29").Value = -40
$ws.Range("N30").Value = -40
